$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1719745222929936
$ws.Range("C2").Value = 0.6178343949044586
$ws.Range("J2").Value = 0.02547770700636943
$ws.Range("P2").Value = 0.1305732484076433
$ws.Range("S2").Value = 0.05414012738853503
$ws.Range("B3").Value = 0.004950495049504951
$ws.Range("C3").Value = 0.01485148514851485
$ws.Range("J3").Value = 0.04455445544554455
$ws.Range("P3").Value = 0.806930693069307
$ws.Range("S3").Value = 0.1287128712871287
$ws.Range("J4").Value = 0.08108108108108109
$ws.Range("P4").Value = 0.7432432432432432
$ws.Range("S4").Value = 0.1756756756756757
$ws.Range("B6").Value = 0.04641350210970464
$ws.Range("D6").Value = 0.02953586497890295
$ws.Range("F6").Value = 0.04641350210970464
$ws.Range("J6").Value = 0.2658227848101266
$ws.Range("O6").Value = 0.01687763713080169
$ws.Range("Q6").Value = 0.189873417721519
$ws.Range("R6").Value = 0.09282700421940929
$ws.Range("S6").Value = 0.3122362869198312
$ws.Range("B7").Value = 0.06930693069306931
$ws.Range("D7").Value = 0.0594059405940594
$ws.Range("F7").Value = 0.07425742574257425
$ws.Range("J7").Value = 0.1485148514851485
$ws.Range("O7").Value = 0.02475247524752475
$ws.Range("Q7").Value = 0.2227722772277228
$ws.Range("R7").Value = 0.07920792079207921
$ws.Range("S7").Value = 0.3217821782178218
$ws.Range("B8").Value = 0.08870967741935484
$ws.Range("D8").Value = 0.02016129032258064
$ws.Range("F8").Value = 0.07258064516129033
$ws.Range("J8").Value = 0.1370967741935484
$ws.Range("O8").Value = 0.02217741935483871
$ws.Range("Q8").Value = 0.1713709677419355
$ws.Range("R8").Value = 0.1088709677419355
$ws.Range("S8").Value = 0.3790322580645161
$ws.Range("B9").Value = 0.09444444444444444
$ws.Range("D9").Value = 0.01666666666666667
$ws.Range("F9").Value = 0.08333333333333333
$ws.Range("J9").Value = 0.1388888888888889
$ws.Range("O9").Value = 0.01666666666666667
$ws.Range("Q9").Value = 0.2611111111111111
$ws.Range("R9").Value = 0.05555555555555555
$ws.Range("B10").Value = 0.1033187226048842
$ws.Range("D10").Value = 0.02629931120851597
$ws.Range("F10").Value = 0.06073888541014402
$ws.Range("J10").Value = 0.1258609893550407
$ws.Range("O10").Value = 0.0100187852222918
$ws.Range("Q10").Value = 0.2298058860363181
$ws.Range("R10").Value = 0.1020663744520977
$ws.Range("S10").Value = 0.3418910457107076
$ws.Range("G11").Value = 0.1471571906354515
$ws.Range("J11").Value = 0.08695652173913043
$ws.Range("K11").Value = 0.2107023411371237
$ws.Range("L11").Value = 0.5451505016722408
$ws.Range("S11").Value = 0.01003344481605351
$ws.Range("G12").Value = 0.7294117647058823
$ws.Range("J12").Value = 0.1764705882352941
$ws.Range("K12").Value = 0.005882352941176471
$ws.Range("L12").Value = 0.05294117647058823
$ws.Range("S12").Value = 0.03529411764705882
$ws.Range("G13").Value = 0.74
$ws.Range("J13").Value = 0.24
$ws.Range("S13").Value = 0.02
$ws.Range("F15").Value = 0.02290076335877863
$ws.Range("H15").Value = 0.1603053435114504
$ws.Range("I15").Value = 0.05725190839694656
$ws.Range("J15").Value = 0.4465648854961832
$ws.Range("K15").Value = 0.04961832061068702
$ws.Range("M15").Value = 0.007633587786259542
$ws.Range("O15").Value = 0.05343511450381679
$ws.Range("S15").Value = 0.2022900763358779
$ws.Range("F16").Value = 0.03305785123966942
$ws.Range("H16").Value = 0.1900826446280992
$ws.Range("I16").Value = 0.06198347107438017
$ws.Range("J16").Value = 0.4545454545454545
$ws.Range("K16").Value = 0.09090909090909091
$ws.Range("M16").Value = 0.02479338842975207
$ws.Range("N16").Value = 0.004132231404958678
$ws.Range("O16").Value = 0.05785123966942149
$ws.Range("S16").Value = 0.08264462809917356
$ws.Range("F17").Value = 0.0135593220338983
$ws.Range("H17").Value = 0.176271186440678
$ws.Range("I17").Value = 0.08135593220338982
$ws.Range("J17").Value = 0.4542372881355932
$ws.Range("K17").Value = 0.09152542372881356
$ws.Range("M17").Value = 0.01864406779661017
$ws.Range("O17").Value = 0.0576271186440678
$ws.Range("S17").Value = 0.1067796610169491
$ws.Range("F18").Value = 0.0193050193050193
$ws.Range("H18").Value = 0.1891891891891892
$ws.Range("I18").Value = 0.07722007722007722
$ws.Range("J18").Value = 0.5057915057915058
$ws.Range("K18").Value = 0.08108108108108109
$ws.Range("M18").Value = 0.007722007722007722
$ws.Range("O18").Value = 0.05019305019305019
$ws.Range("S18").Value = 0.0694980694980695
$ws.Range("F19").Value = 0.01174628034455756
$ws.Range("H19").Value = 0.2028191072826938
$ws.Range("I19").Value = 0.06186374314800313
$ws.Range("J19").Value = 0.4025058731401723
$ws.Range("K19").Value = 0.09240407204385279
$ws.Range("M19").Value = 0.02505873140172279
$ws.Range("N19").Value = 0.001566170712607674
$ws.Range("O19").Value = 0.08613938919342208
$ws.Range("S19").Value = 0.1158966327329679
